$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H38").Value = 282.0909
$ws_ALC.Range("I38").Value = 282.0909
$ws_ALC.Range("K38").Value = 846.2727
$ws_ALC.Range("M38").Value = -474.2727
$ws_ALC.Range("H106").Value = 1023.6667
$ws_ALC.Range("I106").Value = 1023.6667
$ws_ALC.Range("K106").Value = 1023.6667
$ws_ALC.Range("M106").Value = -392.6667
$ws_ALC.Range("H137").Value = 9817.066000000001
$ws_ALC.Range("I137").Value = 14561.333
$ws_ALC.Range("K137").Value = 43683.999
$ws_ALC.Range("M137").Value = -41133.999
$ws_ALC.Range("H138").Value = 6542.46
$ws_ALC.Range("J138").Value = 6627.2314
$ws_ALC.Range("L138").Value = 19881.6942
$ws_ALC.Range("N138").Value = -30161.6942

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 46135.13
$ws_ARM.Range("I32").Value = 36244.344
$ws_ARM.Range("J32").Value = 68243.94
$ws_ARM.Range("K32").Value = 36244.344
$ws_ARM.Range("L32").Value = 68243.94
$ws_ARM.Range("M32").Value = -35957.344
$ws_ARM.Range("N32").Value = -68817.94
$ws_ARM.Range("H61").Value = 8293.5625
$ws_ARM.Range("I61").Value = 8756
$ws_ARM.Range("J61").Value = 8227.5
$ws_ARM.Range("K61").Value = 8756
$ws_ARM.Range("L61").Value = 8227.5
$ws_ARM.Range("M61").Value = -8544
$ws_ARM.Range("N61").Value = -8651.5
$ws_ARM.Range("H74").Value = 8305.777
$ws_ARM.Range("I74").Value = 8854.406000000001
$ws_ARM.Range("K74").Value = 8854.406000000001
$ws_ARM.Range("M74").Value = -7980.406000000001
$ws_ARM.Range("H77").Value = 8305.777
$ws_ARM.Range("I77").Value = 8854.406000000001
$ws_ARM.Range("K77").Value = 44272.03000000001
$ws_ARM.Range("M77").Value = -39904.03000000001
$ws_ARM.Range("H88").Value = 2873.4285
$ws_ARM.Range("I88").Value = 1944.5
$ws_ARM.Range("J88").Value = 3245
$ws_ARM.Range("K88").Value = 1944.5
$ws_ARM.Range("L88").Value = 3245
$ws_ARM.Range("M88").Value = -1538.5
$ws_ARM.Range("N88").Value = -4057
$ws_ARM.Range("H91").Value = 2873.4285
$ws_ARM.Range("I91").Value = 1944.5
$ws_ARM.Range("J91").Value = 3245
$ws_ARM.Range("K91").Value = 1944.5
$ws_ARM.Range("L91").Value = 3245
$ws_ARM.Range("M91").Value = -540.5
$ws_ARM.Range("N91").Value = -6053
$ws_ARM.Range("H135").Value = 0
$ws_ARM.Range("J135").Value = 0
$ws_ARM.Range("L135").Value = 0
$ws_ARM.Range("N135").ClearContents()
$ws_ARM.Range("H136").Value = 8293.5625
$ws_ARM.Range("I136").Value = 8756
$ws_ARM.Range("J136").Value = 8227.5
$ws_ARM.Range("K136").Value = 26268
$ws_ARM.Range("L136").Value = 24682.5
$ws_ARM.Range("M136").Value = -23718
$ws_ARM.Range("N136").Value = -29782.5

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H20").Value = 2889.4285
$ws_BSM.Range("I20").Value = 1908.3636
$ws_BSM.Range("K20").Value = 1908.3636
$ws_BSM.Range("M20").Value = -1661.3636
$ws_BSM.Range("H99").Value = 6529.154
$ws_BSM.Range("I99").Value = 6825.909
$ws_BSM.Range("J99").Value = 4897
$ws_BSM.Range("K99").Value = 6825.909
$ws_BSM.Range("L99").Value = 4897
$ws_BSM.Range("M99").Value = -5327.909
$ws_BSM.Range("N99").Value = -7893
$ws_BSM.Range("H107").Value = 2902.25
$ws_BSM.Range("I107").Value = 3036.3333
$ws_BSM.Range("K107").Value = 3036.3333
$ws_BSM.Range("M107").Value = -1116.3333
$ws_BSM.Range("H134").Value = 3143.7144
$ws_BSM.Range("I134").Value = 3143.7144
$ws_BSM.Range("J134").Value = 0
$ws_BSM.Range("K134").Value = 9431.143199999999
$ws_BSM.Range("L134").Value = 0
$ws_BSM.Range("M134").Value = -6896.143199999999
$ws_BSM.Range("N134").ClearContents()
$ws_BSM.Range("H141").Value = 79999
$ws_BSM.Range("J141").Value = 79999
$ws_BSM.Range("L141").Value = 79999
$ws_BSM.Range("N141").Value = -90359

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 41670820
$ws_CRP.Range("I31").Value = 3465.7
$ws_CRP.Range("K31").Value = 3465.7
$ws_CRP.Range("M31").Value = -3170.7
$ws_CRP.Range("H34").Value = 41670820
$ws_CRP.Range("I34").Value = 3465.7
$ws_CRP.Range("K34").Value = 3465.7
$ws_CRP.Range("M34").Value = -3263.7
$ws_CRP.Range("H99").Value = 3999
$ws_CRP.Range("I99").Value = 3999
$ws_CRP.Range("K99").Value = 3999
$ws_CRP.Range("M99").Value = -2501
$ws_CRP.Range("H107").Value = 4552.5
$ws_CRP.Range("I107").Value = 449.14285
$ws_CRP.Range("K107").Value = 449.14285
$ws_CRP.Range("M107").Value = 1470.85715
$ws_CRP.Range("H126").Value = 3999
$ws_CRP.Range("I126").Value = 3999
$ws_CRP.Range("K126").Value = 11997
$ws_CRP.Range("M126").Value = -9527
$ws_CRP.Range("H132").Value = 3698.8914
$ws_CRP.Range("I132").Value = 2418.0908
$ws_CRP.Range("K132").Value = 7254.2724
$ws_CRP.Range("M132").Value = -4724.2724
$ws_CRP.Range("H141").Value = 218882.78
$ws_CRP.Range("J141").Value = 232098.44
$ws_CRP.Range("L141").Value = 232098.44
$ws_CRP.Range("N141").Value = -242458.44

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H121").Value = 60043140
$ws_CUL.Range("J121").Value = 18635796
$ws_CUL.Range("L121").Value = 55907388
$ws_CUL.Range("N121").Value = -55910008
$ws_CUL.Range("H122").Value = 1824.6
$ws_CUL.Range("I122").Value = 1430.3334
$ws_CUL.Range("K122").Value = 12873.0006
$ws_CUL.Range("M122").Value = -10423.0006

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H126").Value = 18370.785
$ws_GSM.Range("I126").Value = 17385.428
$ws_GSM.Range("J126").Value = 19356.143
$ws_GSM.Range("K126").Value = 52156.284
$ws_GSM.Range("L126").Value = 58068.429
$ws_GSM.Range("M126").Value = -49686.284
$ws_GSM.Range("N126").Value = -63008.429
$ws_GSM.Range("H132").Value = 22267.584
$ws_GSM.Range("I132").Value = 25905.342
$ws_GSM.Range("K132").Value = 77716.026
$ws_GSM.Range("M132").Value = -75186.026

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H22").Value = 2125.7307
$ws_LTW.Range("I22").Value = 1653
$ws_LTW.Range("J22").Value = 2598.4614
$ws_LTW.Range("K22").Value = 1653
$ws_LTW.Range("L22").Value = 2598.4614
$ws_LTW.Range("M22").Value = -1358
$ws_LTW.Range("N22").Value = -3188.4614
$ws_LTW.Range("H27").Value = 2125.7307
$ws_LTW.Range("I27").Value = 1653
$ws_LTW.Range("J27").Value = 2598.4614
$ws_LTW.Range("K27").Value = 1653
$ws_LTW.Range("L27").Value = 2598.4614
$ws_LTW.Range("M27").Value = -1546
$ws_LTW.Range("N27").Value = -2812.4614
$ws_LTW.Range("H40").Value = 778883.25
$ws_LTW.Range("I40").Value = 778883.25
$ws_LTW.Range("K40").Value = 778883.25
$ws_LTW.Range("M40").Value = -778747.25
$ws_LTW.Range("H46").Value = 2669
$ws_LTW.Range("I46").Value = 2655.5
$ws_LTW.Range("J46").Value = 2696
$ws_LTW.Range("K46").Value = 2655.5
$ws_LTW.Range("L46").Value = 2696
$ws_LTW.Range("M46").Value = -2467.5
$ws_LTW.Range("N46").Value = -3072
$ws_LTW.Range("H55").Value = 565.55
$ws_LTW.Range("I55").Value = 501
$ws_LTW.Range("J55").Value = 716.1667
$ws_LTW.Range("K55").Value = 501
$ws_LTW.Range("L55").Value = 716.1667
$ws_LTW.Range("M55").Value = -328
$ws_LTW.Range("N55").Value = -1062.1667

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H94").Value = 36766.332
$ws_WVR.Range("J94").Value = 36766.332
$ws_WVR.Range("L94").Value = 36766.332
$ws_WVR.Range("N94").Value = -38568.332
$ws_WVR.Range("H96").Value = 115093.78
$ws_WVR.Range("I96").Value = 254087
$ws_WVR.Range("J96").Value = 3899.2
$ws_WVR.Range("K96").Value = 254087
$ws_WVR.Range("L96").Value = 3899.2
$ws_WVR.Range("M96").Value = -252714
$ws_WVR.Range("N96").Value = -6645.2
$ws_WVR.Range("H122").Value = 1822.5
$ws_WVR.Range("I122").Value = 1802.7778
$ws_WVR.Range("K122").Value = 5408.3334
$ws_WVR.Range("M122").Value = -2958.3334
$ws_WVR.Range("H126").Value = 2666.5
$ws_WVR.Range("I126").Value = 2778.3
$ws_WVR.Range("K126").Value = 8334.900000000001
$ws_WVR.Range("M126").Value = -5864.900000000001
